$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A width (target stored OOXML width="16.42578125";
# this engine quantizes column width to 1/6-character pixel steps, so the
# closest reachable stored width is 16.5 via ColumnWidth 15.666666666666666)
$ws.Columns.Item(1).ColumnWidth = 15.666666666666666

# New values for A1:A33
$values = @(
    0.31112539381262394,
    -0.057944486788748861,
    -0.017464490136667621,
    -0.0079999998824966667,
    -0.0029999999249126219,
    -0.0019999999078148534,
    -0.0099999998234485687,
    -0.0099999998206476981,
    -0.0019999999012494385,
    -0.0019999998994020274,
    -0.0029999998888339263,
    -0.0034999998835245627,
    -0.0034999998835134605,
    -0.0079999998371445002,
    -0.00099999991076327177,
    0.018799110184872347,
    -0.0019999999002626723,
    -0.0039999998790722913,
    -0.032503384232763999,
    -0.0039999999499702454,
    -0.0039999999494870764,
    -0.0039999999490554217,
    -0.0049999999196401745,
    -0.019999999756725728,
    -0.019999999753369302,
    -0.0024999999002321971,
    -0.002499999895203775,
    -0.001999999875729408,
    -0.0069999998075145342,
    -0.059999999249589653,
    -0.0069999997937468805,
    -0.0099999997619111269,
    -0.003999999823859568
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
